$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:J1): strip " Seen Rx" suffix, uppercase NOCLOG ---
$ws.Range("A1").Value = "FFTR"
$ws.Range("B1").Value = "LIGAZID"
$ws.Range("C1").Value = "EMAZID"
$ws.Range("D1").Value = "LIPICON"
$ws.Range("E1").Value = "AGLIP"
$ws.Range("F1").Value = "CIFIBET"
$ws.Range("G1").Value = "AMLEVO"
$ws.Range("H1").Value = "CARDOBIS"
$ws.Range("I1").Value = "RIVAROX"
$ws.Range("J1").Value = "NOCLOG"

# --- Row labels (column A) and data grid (B:J) for rows 2-42 ---
$rowLabels = @(
  "CMT", "CMT10", "CMT11", "CMT12", "CMT13", "CMT14", "CMT15", "CMT16", "CMT20", "CMT21", "CMT22", "CMT23", "CMT24", "CMT25", "CMT26", "CMT30", "CMT31", "CMT32", "CMT33", "CMT34", "CMT35", "CMT36", "CMT40", "CMT41", "CMT42", "CMT43", "CMT44", "CMT45", "CMT46", "CMT50", "CMT51", "CMT52", "CMT53", "CMT54", "CMT55", "CMT60", "CMT61", "CMT62", "CMT63", "CMT64", "CMT65"
)

$dataGrid = @(
  @(61, 32, 9, 4, 14, 27, 31, 1, 22),
  @(13, 4, 1, 1, 3, 5, 7, 0, 10),
  @(3, 0, 1, 0, 0, 3, 1, 0, 5),
  @(9, 4, 0, 1, 2, 0, 0, 0, 1),
  @(1, 0, 0, 0, 0, 0, 0, 0, 3),
  @(0, 0, 0, 0, 1, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 4, 0, 1),
  @(0, 0, 0, 0, 0, 2, 2, 0, 0),
  @(2, 4, 0, 0, 0, 3, 3, 0, 0),
  @(2, 4, 0, 0, 0, 0, 2, 0, 0),
  @(0, 0, 0, 0, 0, 1, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 2, 1, 0, 0),
  @(10, 11, 0, 0, 0, 5, 6, 0, 1),
  @(8, 9, 0, 0, 0, 0, 0, 0, 1),
  @(2, 2, 0, 0, 0, 5, 6, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(10, 6, 3, 0, 0, 1, 1, 0, 6),
  @(4, 1, 0, 0, 0, 0, 0, 0, 1),
  @(1, 3, 1, 0, 0, 1, 0, 0, 5),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(2, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0, 1, 0, 0, 0, 0, 1, 0, 0),
  @(3, 1, 2, 0, 0, 0, 0, 0, 0),
  @(9, 5, 0, 1, 3, 3, 5, 0, 1),
  @(4, 2, 0, 1, 2, 0, 0, 0, 0),
  @(3, 1, 0, 0, 0, 0, 0, 0, 0),
  @(1, 1, 0, 0, 1, 3, 3, 0, 0),
  @(1, 1, 0, 0, 0, 0, 2, 0, 1),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0),
  @(17, 2, 5, 2, 8, 10, 9, 1, 4),
  @(7, 2, 1, 1, 2, 3, 3, 0, 0),
  @(1, 0, 1, 1, 0, 0, 3, 0, 1),
  @(4, 0, 0, 0, 0, 1, 0, 0, 2),
  @(5, 0, 0, 0, 6, 4, 3, 0, 0),
  @(0, 0, 3, 0, 0, 2, 0, 1, 1)
)

for ($i = 0; $i -lt $rowLabels.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $rowLabels[$i]
  $rowVals = $dataGrid[$i]
  for ($c = 0; $c -lt $rowVals.Count; $c++) {
    $ws.Cells.Item($r, $c + 2).Value = $rowVals[$c]
  }
}
